# "Show error when upload Comment"
#
# The sheet contains a student roster (columns A=Students ID, B=Students Name)
# with a couple of stray "test"/"tse" comment cells left over in columns C/D.
# This edit removes those stray comment cells and instead leaves behind a
# handful of garbled/placeholder values ("asdf", "asd", "d", "adf", ...) in
# several rows - reproducing what the user saw while reproducing the
# "Show error when upload Comment" bug - plus a couple of unrelated
# ID/name corrections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: drop the leftover "test" comment in column C.
$ws.Range("C5").ClearContents()

# Row 6: drop the leftover "test" comment in column D.
$ws.Range("D6").ClearContents()

# New comment text "asdf" typed into rows 13 and 17, column D.
$ws.Range("D13").Value = "asdf"
$ws.Range("D17").Value = "asdf"

# Row 7: the student ID got clobbered with "asd", and a copy of it leaked
# into column D as well.
$ws.Range("A7").Value = "asd"
$ws.Range("D7").Value = "asd"

# Row 10: the student ID here also got clobbered with "asd".
$ws.Range("A10").Value = "asd"

# Row 17: the student name got clobbered with "d".
$ws.Range("B17").Value = "d"

# Row 13: the student name was corrected/retyped.
$ws.Range("B13").Value = "Leartrat Tangvonglearta"

# Row 8: the student ID got clobbered with "adf", and the stray "tse"
# comment in D is removed.
$ws.Range("A8").Value = "adf"
$ws.Range("D8").ClearContents()

# Row 15: the student ID became a bare number (123) instead of text.
$ws.Range("A15").Value = 123

# The author ended up with the cursor sitting on C7 when they saved.
$ws.Range("C7").Select()
